# Set BAEPAbCiPC to 1 for all fuels
$wb = $excel.ActiveWorkbook

$dataSheet = $wb.Worksheets.Item("BAEPAbCiPC")

# Actual fuels (as opposed to energy carriers already =1, or "not a fuel"/"NOT USED" rows)
# get switched on from 0 to 1.
$dataSheet.Range("B3").Value = 1   # hard coal
$dataSheet.Range("B4").Value = 1   # natural gas
$dataSheet.Range("B9:B14").Value = 1  # biomass, petroleum gasoline/diesel, biofuel gasoline/diesel, jet fuel/kerosene
$dataSheet.Range("B17:B20").Value = 1  # lignite, crude oil, heavy fuel oil, LPG propane or butane

# Remove the explanatory paragraph about the U.S. model default carve-out on the About
# sheet, which no longer applies now that all fuels pass through cost changes.
$aboutSheet = $wb.Worksheets.Item("About")
$aboutSheet.Rows("15:18").Delete()
